# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# The report is extended with a new "periodo mora" (2509) for the same
# three workers that already appear for periods 2507 and 2508, and the
# summary totals (VALOR MORA / Cant. Periodos) are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for 3 more data rows (one per worker) right after the
#        existing block of period "2508" rows (16-21), before the blank
#        spacer rows that precede the signature block.
$ws.Rows("22:24").Insert()

# --- 2. Fix up the borders/format of the table now that it has grown:
#        row 21 used to be the last row (with the heavier bottom border);
#        it becomes a normal interior row, and the new row 24 becomes the
#        new last row.
$ws.Range("B21:J21").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B23:J23").PasteSpecial(-4122)   # xlPasteFormats

# --- 3. Fill in the new "2509" rows, mirroring the existing
#        2507 / 2508 blocks (same workers, same Valor Mora / Salario
#        Basico values).
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73131988"
$ws.Range("D22").Value = "DANIEL MORELO MORELO NAVARRO"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "33335743"
$ws.Range("D23").Value = "JOSEFA M GARRIDO CASSIANI"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "19890385"
$ws.Range("D24").Value = "EULOGIO ANTONIO ARRIETA VILLALBA"
$ws.Range("E24").Value = "2509"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

# --- 4. Update the summary header: one more period means 3 periods
#        instead of 2, and the total "Valor Mora" grows from
#        341640 to 512460 (56940 * 3 workers * 3 periods).
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = 512460

Write-Output "edit applied"
